$d = $word.ActiveDocument

# --- Update the date in the first-page header (rId12 / header3.xml) ---
$sections = $d.Sections
$firstSection = $sections.Item(1)
$headers = $firstSection.Headers
$firstPageHeader = $headers.Item(2)
$firstPageHeader.Range.Find.Execute("2023-09-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-15", 2) | Out-Null

# --- Append the new "Knärot" section at the end of the document body ---
# Phase 1: create every paragraph (style + plain text) first, with no
# character-formatting changes interleaved, since applying Font.Italic to a
# range bleeds into formatting used by paragraphs inserted afterwards.
$anchor = $d.Paragraphs.Last

# Paragraph 1
$r = $anchor.Range
$r.InsertParagraphAfter()
$p1 = $d.Paragraphs.Last
$p1.Style = 'Heading 1'
$p1.Range.Text = 'Knärot – ekologi samt krav på livsmiljön'
$anchor = $p1

# Paragraph 2
$r = $anchor.Range
$r.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$p2.Style = 'Normal'
$p2.Range.Text = 'Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021).'
$anchor = $p2

# Paragraph 3
$r = $anchor.Range
$r.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last
$p3.Style = 'Normal'
$p3.Range.Text = 'Samuel Johnsons doktorsavhandling “Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“ (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: “Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” Vidare “More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”'
$anchor = $p3

# Paragraph 4
$r = $anchor.Range
$r.InsertParagraphAfter()
$p4 = $d.Paragraphs.Last
$p4.Style = 'Normal'
$p4.Range.Text = 'Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: “In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”'
$anchor = $p4

# Paragraph 5
$r = $anchor.Range
$r.InsertParagraphAfter()
$p5 = $d.Paragraphs.Last
$p5.Style = 'Normal'
$p5.Range.Text = 'En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022).'
$anchor = $p5

# Paragraph 6
$r = $anchor.Range
$r.InsertParagraphAfter()
$p6 = $d.Paragraphs.Last
$p6.Style = 'Normal'
$p6.Range.Text = 'Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022).'
$anchor = $p6

# Paragraph 7
$r = $anchor.Range
$r.InsertParagraphAfter()
$p7 = $d.Paragraphs.Last
$p7.Style = 'Heading 2'
$p7.Range.Text = 'Referenser - knärot'
$anchor = $p7

# Paragraph 8
$r = $anchor.Range
$r.InsertParagraphAfter()
$p8 = $d.Paragraphs.Last
$p8.Style = 'Normal'
$p8.Range.Text = 'de Graaf M & Roberts M.R., 2009. Short-term response of the herbaceous layer within leave patches after harvest. Forest Ecology and Management 257, 1014-1025'
$anchor = $p8

# Paragraph 9
$r = $anchor.Range
$r.InsertParagraphAfter()
$p9 = $d.Paragraphs.Last
$p9.Style = 'Normal'
$p9.Range.Text = 'Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. Ecological Applications, 22, 2049-2064 '
$anchor = $p9

# Paragraph 10
$r = $anchor.Range
$r.InsertParagraphAfter()
$p10 = $d.Paragraphs.Last
$p10.Style = 'Normal'
$p10.Range.Text = 'Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. Interactive effects of drought and edge exposure on old-growth forest understory species. Landscape Ecology, 37, sid 1839-1853'
$anchor = $p10

# Paragraph 11
$r = $anchor.Range
$r.InsertParagraphAfter()
$p11 = $d.Paragraphs.Last
$p11.Style = 'Normal'
$p11.Range.Text = 'Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. Biological legacies buffer local species extinction after logging. Journal of Applied Ecology. 51, 53-62.'
$anchor = $p11

# Paragraph 12
$r = $anchor.Range
$r.InsertParagraphAfter()
$p12 = $d.Paragraphs.Last
$p12.Style = 'Normal'
$p12.Range.Text = 'Skogsstyrelsen, 2022. Vägledning för hänsyn till knärot. https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/'
$anchor = $p12

# Paragraph 13
$r = $anchor.Range
$r.InsertParagraphAfter()
$p13 = $d.Paragraphs.Last
$p13.Style = 'Normal'
$p13.Range.Text = 'SLU Artdatabanken, 2021. Artfaktablad. Naturvård – artfakta. SLU Artdatabanken, Uppsala '
$anchor = $p13

# Phase 2: now that all paragraphs/structure exist, go back and italicize
# the quoted excerpts / article titles within each paragraph.

# Paragraph 3 italics
$base3 = $p3.Range.Start
$sub = $d.Range($base3 + 34, $base3 + 116)
$sub.Font.Italic = $true
$sub = $d.Range($base3 + 278, $base3 + 483)
$sub.Font.Italic = $true
$sub = $d.Range($base3 + 490, $base3 + 608)
$sub.Font.Italic = $true

# Paragraph 4 italics
$base4 = $p4.Range.Start
$sub = $d.Range($base4 + 205, $base4 + 1070)
$sub.Font.Italic = $true

# Paragraph 8 italics
$base8 = $p8.Range.Start
$sub = $d.Range($base8 + 33, $base8 + 113)
$sub.Font.Italic = $true

# Paragraph 9 italics
$base9 = $p9.Range.Start
$sub = $d.Range($base9 + 62, $base9 + 176)
$sub.Font.Italic = $true

# Paragraph 10 italics
$base10 = $p10.Range.Start
$sub = $d.Range($base10 + 117, $base10 + 207)
$sub.Font.Italic = $true

# Paragraph 11 italics
$base11 = $p11.Range.Start
$sub = $d.Range($base11 + 54, $base11 + 121)
$sub.Font.Italic = $true

# Paragraph 12 italics
$base12 = $p12.Range.Start
$sub = $d.Range($base12 + 22, $base12 + 57)
$sub.Font.Italic = $true

# Paragraph 13 italics
$base13 = $p13.Range.Start
$sub = $d.Range($base13 + 25, $base13 + 61)
$sub.Font.Italic = $true
